# Adds a new weekly price record for "Feria Lagunitas de Puerto Montt" (Espinaca)
# at row 28, shifting all subsequent rows (old 28-64) down by one (new 29-65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 28; this pushes the existing
# rows 28..64 down to 29..65 and grows the used range to A1:R65.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly record.
$ws.Cells.Item(28, 1).Value = 4
$ws.Cells.Item(28, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(28, 3).Value = "Los Lagos"
$ws.Cells.Item(28, 4).Value = 45079
$ws.Cells.Item(28, 5).Value = 10
$ws.Cells.Item(28, 6).Value = 100112012
$ws.Cells.Item(28, 7).Value = "Espinaca"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 25
$ws.Cells.Item(28, 11).Value = 12000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 12000
$ws.Cells.Item(28, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 1200
$ws.Cells.Item(28, 17).Value = 10
$ws.Cells.Item(28, 18).Value = "Hortaliza"
